# Rebrand the "Product Development" multi-industry budget template back to
# "Artificial Intelligence and Machine Learning", restoring the original
# wording (and fixing the mangled "TRProductNING" -> "TRAINING" typo that a
# naive "AI" -> "Product" find/replace had introduced along the way).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Instructions & User Guide
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Instructions & User Guide")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning Comprehensive Budget - User Guide & Instructions"
$ws.Range("A56").Value = "📋 ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING PROJECT OVERVIEW"
$ws.Range("B59").Value = "Data Scientists, ML Engineers, AI Architects, DevOps Engineers..."

# ---------------------------------------------------------------------
# Sheet: Budget Summary
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Budget Summary")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Executive Budget Summary"

# ---------------------------------------------------------------------
# Sheet: Resources
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resources")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Resources Budget"
$ws.Range("A5").Value = "ML Engineers"
$ws.Range("A6").Value = "AI Architects"

# ---------------------------------------------------------------------
# Sheet: Logistics
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Logistics")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Logistics Budget"

# ---------------------------------------------------------------------
# Sheet: Technology
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Technology")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Technology Budget"
$ws.Range("A5").Value = "ML Platform Licenses"

# ---------------------------------------------------------------------
# Sheet: Training
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Training")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Training Budget"
$ws.Range("A4").Value = "AI/ML Certification Programs"
$ws.Range("A10").Value = "TOTAL TRAINING"

# ---------------------------------------------------------------------
# Sheet: Contingency
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contingency")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Contingency Budget"

# ---------------------------------------------------------------------
# Sheet: Timeline
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Timeline")
$ws.Range("A1").Value = "Artificial Intelligence and Machine Learning - Budget Timeline"
